# OpenTBS 1.7.0-beta, can adjust pictures
#
# 1) The sentence that used to be split across two runs around a stray
#    "_GoBack" bookmark ("...automatically mer" | bookmark | "ged in
#    headers and footers.") becomes one contiguous run, and that old
#    bookmark goes away.
# 2) The "changepic" tag gains a new ";adjust" option, typed right before
#    the closing "]"; Word drops its "_GoBack" last-edit bookmark at the
#    new insertion point (immediately before the closing "]").

$d = $word.ActiveDocument

# --- Change 1: merge the split sentence, dropping the old bookmark ----
$mergedFound = $d.Content.Find.Execute( `
    "automatically merged in headers and footers.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "automatically merged in headers and footers.", 2)
Write-Output ("Change 1 (merge sentence) applied: " + $mergedFound)

# --- Change 2: insert ";adjust" before the closing "]" and move the ---
# --- "_GoBack" bookmark to sit right after the inserted text ----------
$rng = $d.Content
$found = $rng.Find.Execute( `
    "[b.number;ope=changepic;from=pic_[val].png;default=current]", `
    $false, $true, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # $rng now spans the whole matched field; its last character is "]".
    $closingBracket = $rng.End - 1

    $insertPoint = $d.Range($closingBracket, $closingBracket)
    $insertPoint.InsertBefore(";adjust")

    # Force the freshly typed text to live in its own run (as real Word
    # does for an in-place insertion) by touching a character format and
    # reverting it immediately.
    $newTextStart = $closingBracket
    $newTextEnd = $closingBracket + 7
    $newTextRange = $d.Range($newTextStart, $newTextEnd)
    $newTextRange.Bold = 1
    $newTextRange.Bold = 0

    # Word leaves its "_GoBack" bookmark at the last edited spot, which
    # also removes any pre-existing "_GoBack" bookmark elsewhere.
    $bmRange = $d.Range($newTextEnd, $newTextEnd)
    $d.Bookmarks.Add("_GoBack", $bmRange)
    Write-Output "Change 2 (;adjust + _GoBack) applied"
} else {
    Write-Output "Change 2 NOT applied: changepic tag text not found"
}
